{"js": "// Remove the trailing \"Ver no Jupiter ...\" / \"\u00a9 2020 ...\" footer block\n// (and the blank paragraph that separated it from the bibliography text),\n// while keeping the blank paragraph / page-break paragraph that originally\n// followed the footer block.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst target = \"Ver no Jupiter Salvar em pdf Salvar em docx\";\nconst copyright = \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\";\n\n// Locate the \"Ver no Jupiter...\" paragraph; the paragraph immediately\n// before it is the blank separator paragraph to remove as well, and the\n// paragraph immediately after it must be the copyright line.\nlet targetIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === target) {\n    targetIndex = i;\n    break;\n  }\n}\n\nif (targetIndex === -1) {\n  throw new Error('Could not find paragraph \"Ver no Jupiter Salvar em pdf Salvar em docx\"');\n}\n\nconst copyrightIndex = targetIndex + 1;\nif (paragraphs.items[copyrightIndex].text !== copyright) {\n  throw new Error(\"Unexpected document structure: copyright paragraph not found right after the target paragraph\");\n}\n\nconst blankIndex = targetIndex - 1;\nif (paragraphs.items[blankIndex].text !== \"\") {\n  throw new Error(\"Unexpected document structure: blank separator paragraph not found right before the target paragraph\");\n}\n\n// Delete from the bottom up so earlier indices stay valid.\nparagraphs.items[copyrightIndex].delete();\nparagraphs.items[targetIndex].delete();\nparagraphs.items[blankIndex].delete();\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" / \"(c) 2020 ...\" footer block\n# (and the blank paragraph that separated it from the bibliography text),\n# while keeping the blank paragraph / page-break paragraph that originally\n# followed the footer block.\n\n$d = $word.ActiveDocument\n\n# Locate the \"Ver no Jupiter...\" paragraph via Find (ASCII needle is safe\n# across the COM text bridge); then map the found range back to its\n# 1-based Paragraphs index by matching start offsets.\n$rng = $d.Content\n$found = $rng.Find.Execute(\"Ver no Jupiter Salvar em pdf Salvar em docx\")\nif (-not $found) {\n    throw \"Could not find paragraph 'Ver no Jupiter Salvar em pdf Salvar em docx'\"\n}\n$targetStart = $rng.Start\n\n$count = $d.Paragraphs.Count\n$targetIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Start -eq $targetStart) {\n        $targetIndex = $i\n        break\n    }\n}\nif ($targetIndex -eq -1) {\n    throw \"Could not map the found text back to a paragraph index\"\n}\n\n# Delete highest index first so the lower indices stay valid:\n#   targetIndex + 1 -> the copyright paragraph\n#   targetIndex      -> \"Ver no Jupiter ...\" paragraph\n#   targetIndex - 1  -> blank separator paragraph\n$d.Paragraphs.Item($targetIndex + 1).Range.Delete()\n$d.Paragraphs.Item($targetIndex).Range.Delete()\n$d.Paragraphs.Item($targetIndex - 1).Range.Delete()\n"}
